# 2022 release build tracker - Synthesis Status update
# Records results of TWIST clonal synthesis run in the "Synthesis Status"
# column (J) of the "Build Plan" sheet for rows that previously had no
# status recorded, plus a few rows whose previously recorded status (all
# "Order placed") is now updated with the actual synthesis result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Build Plan")

# "Failed synthesis" (56 rows)
$failedSynthesisRows = @(
    "J12","J13","J14","J15","J16","J17","J18","J19","J20","J22","J23","J24","J26","J28","J29","J30","J65","J68","J97","J101",
    "J102","J108","J109","J117","J118","J119","J120","J121","J124","J128","J129","J131","J138","J142","J146","J150","J153","J155","J167","J171",
    "J175","J178","J188","J189","J191","J197","J209","J213","J214","J215","J217","J219","J226","J230","J232","J240"
)
foreach ($cellAddr in $failedSynthesisRows) {
    $ws.Range($cellAddr).Value = "Failed synthesis"
}

# "Available in FreeGenes" (104 rows)
$availableInFreeGenesRows = @(
    "J21","J25","J27","J31","J32","J33","J34","J37","J38","J47","J48","J49","J51","J52","J53","J55","J56","J60","J61","J63",
    "J67","J69","J71","J72","J74","J77","J79","J80","J81","J85","J88","J89","J90","J92","J94","J96","J98","J99","J100","J103",
    "J104","J105","J106","J107","J111","J112","J113","J114","J115","J116","J125","J126","J127","J130","J132","J133","J135","J137","J141","J144",
    "J145","J147","J148","J149","J151","J152","J154","J157","J161","J163","J164","J165","J166","J168","J172","J173","J177","J180","J182","J192",
    "J193","J196","J206","J207","J210","J211","J212","J218","J220","J221","J224","J225","J227","J228","J229","J233","J234","J235","J236","J237",
    "J238","J239","J241","J242"
)
foreach ($cellAddr in $availableInFreeGenesRows) {
    $ws.Range($cellAddr).Value = "Available in FreeGenes"
}

# "Order placed" (10 rows)
$orderPlacedRows = @(
    "J158","J159","J162","J184","J186","J190","J208","J222","J223","J231"
)
foreach ($cellAddr in $orderPlacedRows) {
    $ws.Range($cellAddr).Value = "Order placed"
}

# Restore the navigation/selection state recorded for the sheet: the
# viewport top-left cell is E6 and the active selection is I25.
$ws.Activate()
$ws.Range("I25").Select()
try {
    $excel.ActiveWindow.ScrollRow = 6
    $excel.ActiveWindow.ScrollColumn = 5
} catch {
    # Scroll position APIs may not be fully supported; selection is the
    # important, persisted part of the view state.
}
